# Update to framework 7.1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Add the missing "CREATE/MODIFY" value in column A of row 5, matching the
# same value/style already used in A2:A4 (same font as the cells above).
$ws.Range("A5").Value = "CREATE/MODIFY"
$ws.Range("A5").Font.Name = $ws.Range("A2").Font.Name
$ws.Range("A5").Font.Size = $ws.Range("A2").Font.Size
$ws.Range("A5").Font.ColorIndex = $ws.Range("A2").Font.ColorIndex

# Update the active sheet view: drop the frozen/scrolled topLeftCell and
# move the current selection from F12 to B9.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select() | Out-Null
